$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '62.463.72'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.25%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.192.91'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.45%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '593.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.78'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.64%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.189.83'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.507'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.95%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.144'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.94%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.32'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.82%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.454'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -3.76%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000238'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.96%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '33.56'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.85%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.724.55'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -3.28%  '
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.197.55'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '62.577.74'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -2.22%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.02%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '462.86'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.03%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.06'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.712'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.14%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.68'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -4.34%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.51'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.23%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.97'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.21%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('D26').Style = 'Normal'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.71'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.60%  '
$ws.Range('E28').Value = '  +0.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.93'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -4.31%  '
$ws.Range('E30').Value = '  -5.02%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.07'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.87%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '27.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.59%  '
$ws.Range('E33').Value = '  -4.40%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.44'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.55%  '
$ws.Range('E35').Value = '  -4.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.86'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.29%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '51.55'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.38%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0₃0693'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -9.47%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0392'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.90%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.009.10'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.75%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '415.06'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.115'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +4.87%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.10'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.63'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -6.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.253'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.18%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.17'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.10%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '36.07'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.27%  '
$ws.Range('E48').Value = '  -0.07%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.91'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -2.22%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '124.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.83%  '
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.30'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.65%  '
